$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.816.68"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.646.27"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.17%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.56%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.34"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.503"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.58%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.54%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0627"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.23"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.49%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.867.64"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.69%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.634.51"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.19"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.77%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.528"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.77"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.75%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.815.73"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0739"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.85%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "214.70"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.75%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.70%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.43"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +13.76%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.30"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.80%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.38"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.31"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.17%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.119"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.78%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.11"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.73"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.95%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0516"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.34%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.54%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.50%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.288.26"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.10%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.24%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.28%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.61%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.828"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.811"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.24%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.23%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.793.81"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.77"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.81%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.23"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.04%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.58%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.90%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.70"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0980"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.59%  "
